# Updated symbol list on Tue Dec 27 21:43:17 UTC 2022 with GitHub Actions
#
# Applies the latest coinranking.com price/volume refresh to Sheet1:
#  - numeric-looking "Price" (column D) cells are re-written as literal text
#    (matching the existing inlineStr layout), using a leading apostrophe so
#    the interop layer doesn't silently coerce them into numeric cells, then
#    the style is reset back to Normal so no stray NumberFormat/quotePrefix
#    style sticks around on the cell.
#  - two pairs of rows (6/7 and 41/42/43) had their coin ranking reordered;
#    Coin / Link / Price / Volume(1h) are updated together for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$cellRef,
        [string]$text
    )
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# --- simple price refreshes -------------------------------------------------
Set-TextValue "D2" "245.27"
Set-TextValue "D3" "24.00"
Set-TextValue "D4" "5.358"

# --- row 6 / row 7 swap (KuCoinToken <-> GateToken) -------------------------
Set-TextValue "B6" "GateToken"
Set-TextValue "C6" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D6" "3.369"
Set-TextValue "E6" "5GateTokenGT"

Set-TextValue "B7" "KuCoinToken"
Set-TextValue "C7" "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue "D7" "6.465"
Set-TextValue "E7" "6KuCoinTokenKCS"

# --- more price refreshes ----------------------------------------------------
Set-TextValue "D8" "0.8103"
Set-TextValue "D9" "0.9178"
Set-TextValue "D10" "0.1404"
Set-TextValue "D11" "0.07382"
Set-TextValue "D12" "0.03173"
Set-TextValue "D13" "0.03067"
Set-TextValue "D14" "0.09361"
Set-TextValue "D15" "3.848"
Set-TextValue "D16" "0.001560"
Set-TextValue "D17" "0.04695"
Set-TextValue "D18" "0.0005981"
Set-TextValue "E18" "17OneONEWorstin24h"
Set-TextValue "D19" "0.006074"
Set-TextValue "D20" "0.001247"
Set-TextValue "D21" "0.004686"
Set-TextValue "D22" "0.00008798"
Set-TextValue "D23" "3.593"
Set-TextValue "D28" "0.0002350"
Set-TextValue "D40" "0.03840"

# --- rows 41 / 42 / 43 three-way rotation (BKEXToken -> KickToken,
#     CEJI -> BKEXToken, KickToken -> CEJI) -----------------------------------
Set-TextValue "B41" "KickToken"
Set-TextValue "C41" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006298"
Set-TextValue "E41" "40KickTokenKICK"

Set-TextValue "B42" "BKEXToken"
Set-TextValue "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1064"
Set-TextValue "E42" "41BKEXTokenBKK"

Set-TextValue "B43" "CEJI"
Set-TextValue "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003199"
Set-TextValue "E43" "42CEJICEJI"

# --- trailing price refreshes ------------------------------------------------
Set-TextValue "D44" "0.009028"
Set-TextValue "D47" "0.6871"
Set-TextValue "D48" "0.001855"
